# Auto-generated edit script: updates crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.032.85'
$ws.Range("E2").Value = '  -1.07%  '
$ws.Range("D3").Value = '3.049.03'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.18%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '386.88'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.15%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '101.86'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.49%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.534'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.10%  '
$ws.Range("E8").Value = '  +0.08%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.577'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.12%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '36.50'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -0.54%  '
$ws.Range("E11").Value = '  +0.03%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0846'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.85%  '
$ws.Range("D13").Value = '3.549.65'
$ws.Range("E13").Value = '  +1.52%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '18.23'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.93%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("D16").Value = '3.018.54'
$ws.Range("E16").Value = '  -0.13%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.981'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.93%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '10.60'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").Value = '51.075.58'
$ws.Range("E19").Value = '  -0.95%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '3.19'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0954'
$ws.Range("E21").Value = '  -0.92%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '12.20'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.86%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '69.51'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.59%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '263.33'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("E25").Value = '  -0.84%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '7.86'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -5.34%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '26.87'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.70%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '7.19'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -4.87%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -4.94%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.104'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.28%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '10.38'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.26%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '35.20'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +4.14%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0467'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +4.07%  '
$ws.Range("E35").Value = '  -0.10%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '50.02'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  +1.85%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.286'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -1.62%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '129.30'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +1.96%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.82'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.47%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.114'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.25%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.79'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '16.28'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.54%  '
$ws.Range("E45").Value = '  -2.90%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '21.53'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '2.48'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +3.70%  '
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '2.054.66'
$ws.Range("E49").Value = '  +1.45%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '9.37'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +11.19%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.907'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +15.50%  '
